$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45923
$ws.Range("B2").Value = 5046.48585762747
$ws.Range("C2").Value = 4954.73008201507
$ws.Range("D2").Value = 6432
$ws.Range("E2").Value = 6523.4375
$ws.Range("F2").Value = -0.0132614838497223

$ws.Range("A3").Value = 45924
$ws.Range("B3").Value = 5032.32142458527
$ws.Range("C3").Value = 5102.01014392451
$ws.Range("D3").Value = 2952
$ws.Range("E3").Value = 6504.363296
$ws.Range("F3").Value = 150.918833972468

$ws.Range("A4").Value = 45925
$ws.Range("B4").Value = 5051.33999310527
$ws.Range("C4").Value = 5573.73971489882
$ws.Range("D4").Value = 2952
$ws.Range("E4").Value = 6525.782493
$ws.Range("F4").Value = 170.674258949731

$ws.Range("A5").Value = 45926
$ws.Range("B5").Value = 5096.26951046918
$ws.Range("C5").Value = 4957.55197502764
$ws.Range("D5").Value = 2952
$ws.Range("E5").Value = 6576.564692
$ws.Range("F5").Value = 145.243631523269

$ws.Range("A6").Value = 45927
$ws.Range("B6").Value = 1315.89070720157
$ws.Range("C6").Value = 2994.97829561339
$ws.Range("D6").Value = 2952
$ws.Range("E6").Value = 2527.530788
$ws.Range("F6").Value = 52.2757656838257

$ws.Range("A7").Value = 45928
$ws.Range("B7").Value = 1089.40583490764
$ws.Range("C7").Value = 2755.76463151152
$ws.Range("D7").Value = 2952
$ws.Range("E7").Value = 2183.81959
$ws.Range("F7").Value = 37.4240994418282

$ws.Range("A8").Value = 45929
$ws.Range("B8").Value = 5457.44872542722
$ws.Range("C8").Value = 5451.20771697766
$ws.Range("D8").Value = 2952
$ws.Range("E8").Value = 7030.796743
$ws.Range("F8").Value = 169.689822272935

$ws.Range("A9").Value = 45930
$ws.Range("B9").Value = 5457.44872542722
$ws.Range("C9").Value = 5651.30360979221
$ws.Range("D9").Value = 2952
$ws.Range("E9").Value = 7030.796743
$ws.Range("F9").Value = 178.027151140208

$ws.Range("A10").Value = 45931
$ws.Range("B10").Value = 4260.05383201679
$ws.Range("C10").Value = 4723.17619978083
$ws.Range("D10").Value = 3692
$ws.Range("E10").Value = 6097.879545
$ws.Range("F10").Value = 119.541746365168

$ws.Range("A11").Value = 45932
$ws.Range("B11").Value = 4260.05383201679
$ws.Range("C11").Value = 4676.09676603563
$ws.Range("D11").Value = 3692
$ws.Range("E11").Value = 6097.879545
$ws.Range("F11").Value = 117.580103292452

$ws.Range("A12").Value = 45933
$ws.Range("B12").Value = 4260.05383201679
$ws.Range("C12").Value = 4055.7375268143
$ws.Range("D12").Value = 3692
$ws.Range("E12").Value = 6097.879545
$ws.Range("F12").Value = 91.7318016582296

$ws.Range("A13").Value = 45934
$ws.Range("B13").Value = 805.573379841386
$ws.Range("C13").Value = 2217.56746708433
$ws.Range("D13").Value = 3692
$ws.Range("E13").Value = 2301.949365
$ws.Range("F13").Value = 0.914310510122637

$ws.Range("A14").Value = 45935
$ws.Range("B14").Value = 709.592596030573
$ws.Range("C14").Value = 2172.18673418899
$ws.Range("D14").Value = 3692
$ws.Range("E14").Value = 2197.537945
$ws.Range("F14").Value = -1.32782986839946

$ws.Range("A15").Value = 45936
$ws.Range("B15").Value = 4367.23135895568
$ws.Range("C15").Value = 4652.85285309916
$ws.Range("D15").Value = 3692
$ws.Range("E15").Value = 6290.949659
$ws.Range("F15").Value = 120.190464714311
